$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.522.66'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '1.958.86'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Formula = "'244.54"
$ws.Range('E5').Value = '  +0.80%  '
$ws.Range('D6').Formula = "'0.622"
$ws.Range('E6').Value = '  -0.42%  '
$ws.Range('D7').Formula = "'58.65"
$ws.Range('E7').Value = '  -1.77%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -1.65%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Formula = "'0.0848"
$ws.Range('E10').Value = '  +5.88%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Formula = "'0.104"
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Formula = "'22.07"
$ws.Range('E12').Value = '  -1.73%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Formula = "'0.832"
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.245.87'
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Formula = "'13.67"
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Formula = "'5.26"
$ws.Range('E16').Value = '  -3.27%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.961.13'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '36.444.06'
$ws.Range('E18').Value = '  +0.83%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Formula = "'70.10"
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0867'
$ws.Range('E20').Value = '  +1.21%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Formula = "'230.27"
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Formula = "'5.07"
$ws.Range('E22').Value = '  -2.85%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Formula = "'1.00"
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Formula = "'2.47"
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Formula = "'2.31"
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Formula = "'9.31"
$ws.Range('E26').Value = '  -5.35%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Formula = "'0.139"
$ws.Range('E27').Value = '  +10.35%  '
$ws.Range('D28').Formula = "'162.21"
$ws.Range('E28').Value = '  +1.15%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Formula = "'19.56"
$ws.Range('E29').Value = '  -1.49%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Formula = "'0.119"
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Formula = "'1.18"
$ws.Range('E31').Value = '  +3.37%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Formula = "'4.72"
$ws.Range('E32').Value = '  -3.24%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Formula = "'0.0635"
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D34').Formula = "'4.31"
$ws.Range('E34').Value = '  -2.21%  '
$ws.Range('B35').Value = 'THORChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D35').Formula = "'6.38"
$ws.Range('E35').Value = '  +1.42%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Formula = "'1.00"
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Formula = "'1.77"
$ws.Range('E37').Value = '  -2.19%  '
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').Formula = "'2.18"
$ws.Range('E38').Value = '  -4.47%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Formula = "'3.06"
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('B40').Value = 'Cronos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D40').Formula = "'0.0989"
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').Formula = "'2.87"
$ws.Range('E41').Value = '  +0.38%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Formula = "'1.18"
$ws.Range('E42').Value = '  -3.48%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Formula = "'0.0211"
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Formula = "'16.18"
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.366.23'
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('B46').Value = 'ARBITRUM'
$ws.Range('C46').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D46').Formula = "'1.04"
$ws.Range('E46').Value = '  -4.20%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Formula = "'88.56"
$ws.Range('E47').Value = '  -4.14%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Formula = "'7.23"
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Formula = "'2.83"
$ws.Range('E49').Value = '  -0.23%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Formula = "'46.06"
$ws.Range('E50').Value = '  +4.26%  '
$ws.Range('B51').Value = 'RocketPoolETH'
$ws.Range('C51').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D51').Value = '2.135.45'
$ws.Range('E51').Value = '  -0.61%  '
